$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 7
$ws.Range("H7").Value = 6475
$ws.Range("J7").Value = 2950
$ws.Range("L7").Value = 2950
$ws.Range("N7").Value = -3174

# Row 10
$ws.Range("H10").Value = 2950
$ws.Range("J10").Value = 2950
$ws.Range("L10").Value = 2950
$ws.Range("N10").Value = -3536

# Row 14
$ws.Range("H14").Value = 6475
$ws.Range("J14").Value = 2950
$ws.Range("L14").Value = 2950
$ws.Range("N14").Value = -3332

# Row 53
$ws.Range("H53").Value = 199.35715
$ws.Range("I53").Value = 89.8
$ws.Range("J53").Value = 260.22223
$ws.Range("K53").Value = 89.8
$ws.Range("L53").Value = 260.22223
$ws.Range("M53").Value = 547.2
$ws.Range("N53").Value = -1534.22223

# Row 137
$ws.Range("H137").Value = 1521.8518
$ws.Range("I137").Value = 969.35
$ws.Range("J137").Value = 3100.4285
$ws.Range("K137").Value = 2908.05
$ws.Range("L137").Value = 9301.2855
$ws.Range("M137").Value = -358.0500000000002
$ws.Range("N137").Value = -14401.2855

$ws = $wb.Worksheets.Item("ARM")
# Row 22
$ws.Range("H22").Value = 9075
$ws.Range("I22").Value = 433.33334
$ws.Range("J22").Value = 35000
$ws.Range("K22").Value = 433.33334
$ws.Range("L22").Value = 35000
$ws.Range("M22").Value = -134.33334
$ws.Range("N22").Value = -35598

# Row 26
$ws.Range("H26").Value = 798.5
$ws.Range("I26").Value = 798.5
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 798.5
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -468.5
$ws.Range("N26").ClearContents()

# Row 32
$ws.Range("H32").Value = 11653.015
$ws.Range("I32").Value = 7931.1055
$ws.Range("J32").Value = 29332.084
$ws.Range("K32").Value = 7931.1055
$ws.Range("L32").Value = 29332.084
$ws.Range("M32").Value = -7644.1055
$ws.Range("N32").Value = -29906.084

# Row 61
$ws.Range("H61").Value = 1750.3429
$ws.Range("I61").Value = 1574.9333
$ws.Range("J61").Value = 2802.8
$ws.Range("K61").Value = 1574.9333
$ws.Range("L61").Value = 2802.8
$ws.Range("M61").Value = -1362.9333
$ws.Range("N61").Value = -3226.8

# Row 74
$ws.Range("H74").Value = 23810784
$ws.Range("I74").Value = 23810784
$ws.Range("K74").Value = 23810784
$ws.Range("M74").Value = -23809910

# Row 77
$ws.Range("H77").Value = 23810784
$ws.Range("I77").Value = 23810784
$ws.Range("K77").Value = 119053920
$ws.Range("M77").Value = -119049552

# Row 102
$ws.Range("H102").Value = 1636.0303
$ws.Range("I102").Value = 1454.7931
$ws.Range("J102").Value = 2950
$ws.Range("K102").Value = 1454.7931
$ws.Range("L102").Value = 2950
$ws.Range("M102").Value = 167.2068999999999
$ws.Range("N102").Value = -6194

# Row 136
$ws.Range("H136").Value = 1750.3429
$ws.Range("I136").Value = 1574.9333
$ws.Range("J136").Value = 2802.8
$ws.Range("K136").Value = 4724.7999
$ws.Range("L136").Value = 8408.400000000001
$ws.Range("M136").Value = -2174.7999
$ws.Range("N136").Value = -13508.4

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 11113518
$ws.Range("I31").Value = 17242344
$ws.Range("J31").Value = 5018.9375
$ws.Range("K31").Value = 17242344
$ws.Range("L31").Value = 5018.9375
$ws.Range("M31").Value = -17242049
$ws.Range("N31").Value = -5608.9375

# Row 34
$ws.Range("H34").Value = 11113518
$ws.Range("I34").Value = 17242344
$ws.Range("J34").Value = 5018.9375
$ws.Range("K34").Value = 17242344
$ws.Range("L34").Value = 5018.9375
$ws.Range("M34").Value = -17242142
$ws.Range("N34").Value = -5422.9375

# Row 58
$ws.Range("H58").Value = 1567.4634
$ws.Range("I58").Value = 835.0769
$ws.Range("J58").Value = 2836.9333
$ws.Range("K58").Value = 835.0769
$ws.Range("L58").Value = 2836.9333
$ws.Range("M58").Value = -632.0769
$ws.Range("N58").Value = -3242.9333

# Row 94
$ws.Range("H94").Value = 3289.1667
$ws.Range("J94").Value = 2973.111
$ws.Range("L94").Value = 2973.111
$ws.Range("N94").Value = -3875.111

# Row 122
$ws.Range("H122").Value = 1592.375
$ws.Range("I122").Value = 1186.5714
$ws.Range("J122").Value = 1810.8846
$ws.Range("K122").Value = 3559.7142
$ws.Range("L122").Value = 5432.6538
$ws.Range("M122").Value = -1109.7142
$ws.Range("N122").Value = -10332.6538

# Row 132
$ws.Range("H132").Value = 1364.25
$ws.Range("I132").Value = 930.7083
$ws.Range("J132").Value = 2664.875
$ws.Range("K132").Value = 2792.1249
$ws.Range("L132").Value = 7994.625
$ws.Range("M132").Value = -262.1248999999998
$ws.Range("N132").Value = -13054.625

# Row 136
$ws.Range("H136").Value = 1567.4634
$ws.Range("I136").Value = 835.0769
$ws.Range("J136").Value = 2836.9333
$ws.Range("K136").Value = 2505.2307
$ws.Range("L136").Value = 8510.7999
$ws.Range("M136").Value = 44.76929999999993
$ws.Range("N136").Value = -13610.7999

# Row 138
$ws.Range("H138").Value = 53000
$ws.Range("J138").Value = 51600
$ws.Range("L138").Value = 51600
$ws.Range("N138").Value = -61880

# Row 140
$ws.Range("H140").Value = 48593.332
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 48593.332
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 48593.332
$ws.Range("M140").ClearContents()
$ws.Range("N140").Value = -58953.332

$ws = $wb.Worksheets.Item("CUL")
# Row 22
$ws.Range("H22").Value = 2000
$ws.Range("J22").Value = 3000
$ws.Range("L22").Value = 9000
$ws.Range("N22").Value = -9338

# Row 27
$ws.Range("H27").Value = 2000
$ws.Range("J27").Value = 3000
$ws.Range("L27").Value = 9000
$ws.Range("N27").Value = -9204

# Row 50
$ws.Range("H50").Value = 293.72726
$ws.Range("I50").Value = 251.66667
$ws.Range("J50").Value = 309.5
$ws.Range("K50").Value = 755.00001
$ws.Range("L50").Value = 928.5
$ws.Range("M50").Value = -274.00001
$ws.Range("N50").Value = -1890.5

# Row 53
$ws.Range("H53").Value = 293.72726
$ws.Range("I53").Value = 251.66667
$ws.Range("J53").Value = 309.5
$ws.Range("K53").Value = 755.00001
$ws.Range("L53").Value = 928.5
$ws.Range("M53").Value = -274.00001
$ws.Range("N53").Value = -1890.5

# Row 98
$ws.Range("H98").Value = 195.42857
$ws.Range("I98").Value = 210
$ws.Range("J98").Value = 184.5
$ws.Range("K98").Value = 630
$ws.Range("L98").Value = 553.5
$ws.Range("M98").Value = 868
$ws.Range("N98").Value = -3549.5

# Row 131
$ws.Range("H131").Value = 84889.72
$ws.Range("I131").Value = 700
$ws.Range("J131").Value = 92210.56
$ws.Range("K131").Value = 2100
$ws.Range("L131").Value = 276631.68
$ws.Range("M131").Value = 2940
$ws.Range("N131").Value = -286711.68

$ws = $wb.Worksheets.Item("GSM")
# Row 17
$ws.Range("H17").Value = 25000
$ws.Range("J17").Value = 25000
$ws.Range("L17").Value = 25000
$ws.Range("N17").Value = -25336

# Row 20
$ws.Range("H20").Value = 10000
$ws.Range("J20").Value = 10000
$ws.Range("L20").Value = 10000
$ws.Range("N20").Value = -10490

# Row 97
$ws.Range("H97").Value = 892.4483
$ws.Range("I97").Value = 726.8182
$ws.Range("J97").Value = 1413
$ws.Range("K97").Value = 726.8182
$ws.Range("L97").Value = 1413
$ws.Range("M97").Value = -230.8182
$ws.Range("N97").Value = -2405

# Row 132
$ws.Range("H132").Value = 5443.7812
$ws.Range("I132").Value = 6590.5
$ws.Range("J132").Value = 3532.5833
$ws.Range("K132").Value = 19771.5
$ws.Range("L132").Value = 10597.7499
$ws.Range("M132").Value = -17241.5
$ws.Range("N132").Value = -15657.7499

$ws = $wb.Worksheets.Item("LTW")
# Row 5
$ws.Range("H5").Value = 9000
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 9000
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 9000
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -9226

# Row 82
$ws.Range("H82").Value = 1547.2609
$ws.Range("I82").Value = 1611.1666
$ws.Range("J82").Value = 1477.5454
$ws.Range("K82").Value = 1611.1666
$ws.Range("L82").Value = 1477.5454
$ws.Range("M82").Value = -1250.1666
$ws.Range("N82").Value = -2199.5454

# Row 85
$ws.Range("H85").Value = 1547.2609
$ws.Range("I85").Value = 1611.1666
$ws.Range("J85").Value = 1477.5454
$ws.Range("K85").Value = 1611.1666
$ws.Range("L85").Value = 1477.5454
$ws.Range("M85").Value = -363.1666
$ws.Range("N85").Value = -3973.5454

# Row 106
$ws.Range("H106").Value = 21092.5
$ws.Range("J106").Value = 21092.5
$ws.Range("L106").Value = 21092.5
$ws.Range("N106").Value = -23616.5

# Row 112
$ws.Range("H112").Value = 16816.883
$ws.Range("J112").Value = 16816.883
$ws.Range("L112").Value = 16816.883
$ws.Range("N112").Value = -19770.883

# Row 132
$ws.Range("H132").Value = 2224.6667
$ws.Range("I132").Value = 1477.5769
$ws.Range("J132").Value = 4999.5713
$ws.Range("K132").Value = 4432.7307
$ws.Range("L132").Value = 14998.7139
$ws.Range("M132").Value = -1902.7307
$ws.Range("N132").Value = -20058.7139

$ws = $wb.Worksheets.Item("WVR")
# Row 24
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()

# Row 105
$ws.Range("H105").Value = 35000
$ws.Range("J105").Value = 35000
$ws.Range("L105").Value = 35000
$ws.Range("N105").Value = -41988
